$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the report date-range header (A4, merged A4:I4)
$ws.Range("A4").Value = "Từ ngày: 01-08-2022 đến ngày: 30-10-2022"

# 2. Insert 5 new data rows right after the last existing item row (row 21),
#    which pushes the old footer rows (22-25, 29) down to (27-30, 34).
$ws.Rows("22:26").Insert()

# 3. Give the newly inserted rows the same bordered look as the other
#    item rows (A6:I21) before filling in their values.
for ($row = 22; $row -le 26; $row++) {
    $ws.Range("A" + $row + ":I" + $row).Borders.LineStyle = 1
}

# 4. Fill in the new item rows.
#    Columns: A=STT, B=Hàng hóa, C=Danh mục, D=ĐVT, E=Mã hàng hóa,
#             F=Giá Bán, G=Số lượng, H=Thành tiền, I=Ghi Chú
$ws.Cells.Item(22, 1).Value = 17
$ws.Cells.Item(22, 2).Value = "Cad điện thoại"
$ws.Cells.Item(22, 3).Value = "TSDTK"
$ws.Cells.Item(22, 4).Value = "cái"
$ws.Cells.Item(22, 5).Value = "B001"
$ws.Cells.Item(22, 6).Value = "10,000"
$ws.Cells.Item(22, 7).Value = 1
$ws.Cells.Item(22, 8).Value = 10000

$ws.Cells.Item(23, 1).Value = 18
$ws.Cells.Item(23, 2).Value = "Cad điện thoại"
$ws.Cells.Item(23, 3).Value = "TSDTK"
$ws.Cells.Item(23, 4).Value = "cái"
$ws.Cells.Item(23, 5).Value = "B004"
$ws.Cells.Item(23, 6).Value = "40,000"
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(23, 8).Value = 40000

$ws.Cells.Item(24, 1).Value = 19
$ws.Cells.Item(24, 2).Value = "Cad điện thoại"
$ws.Cells.Item(24, 3).Value = "TSDTK"
$ws.Cells.Item(24, 4).Value = "cái"
$ws.Cells.Item(24, 5).Value = "B008"
$ws.Cells.Item(24, 6).Value = "80,000"
$ws.Cells.Item(24, 7).Value = 1
$ws.Cells.Item(24, 8).Value = 80000

$ws.Cells.Item(25, 1).Value = 20
$ws.Cells.Item(25, 2).Value = "Cad điện thoại"
$ws.Cells.Item(25, 3).Value = "TSDTK"
$ws.Cells.Item(25, 4).Value = "cái"
$ws.Cells.Item(25, 5).Value = "B008"
$ws.Cells.Item(25, 6).Value = "80,000"
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 80000

$ws.Cells.Item(26, 1).Value = 21
$ws.Cells.Item(26, 2).Value = "Nước yến"
$ws.Cells.Item(26, 3).Value = "DK"
$ws.Cells.Item(26, 4).Value = "Lon"
$ws.Cells.Item(26, 5).Value = "CT5Z"
$ws.Cells.Item(26, 6).Value = "11,000"
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 11000

# 5. Update the totals / signature block that got shifted down to rows 27-30 and 34.
$ws.Cells.Item(27, 8).Value = "726,500"
$ws.Cells.Item(28, 2).Value = "bảy trăm  hai mươi sáu nghìn năm trăm đồng"
$ws.Cells.Item(29, 7).Value = "ngày 03 tháng 10 năm 2022"

# 6. Widen column B (Hàng hóa) from 39 to 50 characters.
$ws.Columns(2).ColumnWidth = 49.35

# 7. Update the selected cell shown when the sheet is opened.
$ws.Range("I26").Select()
